$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accidentes_MEX")

$ws.Range("B2").Value = 331938
$ws.Range("B20").Value = 355281

$wb.Save()
